# Masamune_Profits workbook update (scheduled runner refresh of market-price
# derived columns H:N - currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ),
# LeveProfit(NQ/HQ) - across the ALC/ARM/BSM/CRP/CUL/LTW/WVR leve tables).
# Values below were recomputed upstream; this script just writes the
# refreshed numbers into the matching cells (some rows also gain
# previously-empty LeveProfit cells).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 4199.5  # H69
$ws.Cells.Item(69, 9).Value = 3800  # I69
$ws.Cells.Item(69, 10).Value = 4599  # J69
$ws.Cells.Item(69, 11).Value = 11400  # K69
$ws.Cells.Item(69, 12).Value = 13797  # L69
$ws.Cells.Item(69, 13).Value = -10526  # M69
$ws.Cells.Item(69, 14).Value = -15545  # N69
$ws.Cells.Item(72, 8).Value = 4199.5  # H72
$ws.Cells.Item(72, 9).Value = 3800  # I72
$ws.Cells.Item(72, 10).Value = 4599  # J72
$ws.Cells.Item(72, 11).Value = 34200  # K72
$ws.Cells.Item(72, 12).Value = 41391  # L72
$ws.Cells.Item(72, 13).Value = -29832  # M72
$ws.Cells.Item(72, 14).Value = -50127  # N72
$ws.Cells.Item(116, 8).Value = 6055.5557  # H116
$ws.Cells.Item(116, 9).Value = 3524.75  # I116
$ws.Cells.Item(116, 10).Value = 8080.2  # J116
$ws.Cells.Item(116, 11).Value = 3524.75  # K116
$ws.Cells.Item(116, 12).Value = 8080.2  # L116
$ws.Cells.Item(116, 13).Value = -82.75  # M116
$ws.Cells.Item(116, 14).Value = -14964.2  # N116
$ws.Cells.Item(125, 8).Value = 3112.75  # H125
$ws.Cells.Item(125, 9).Value = 2620.4  # I125
$ws.Cells.Item(125, 10).Value = 3933.3333  # J125
$ws.Cells.Item(125, 11).Value = 23583.6  # K125
$ws.Cells.Item(125, 12).Value = 35399.9997  # L125
$ws.Cells.Item(125, 13).Value = -21123.6  # M125
$ws.Cells.Item(125, 14).Value = -40319.9997  # N125
$ws.Cells.Item(131, 8).Value = 2528  # H131
$ws.Cells.Item(131, 9).Value = 3795  # I131
$ws.Cells.Item(131, 10).Value = 2105.6667  # J131
$ws.Cells.Item(131, 11).Value = 11385  # K131
$ws.Cells.Item(131, 12).Value = 6317.000100000001  # L131
$ws.Cells.Item(131, 13).Value = -6345  # M131
$ws.Cells.Item(131, 14).Value = -16397.0001  # N131
$ws.Cells.Item(132, 8).Value = 20850.375  # H132
$ws.Cells.Item(132, 9).Value = 2845.3555  # I132
$ws.Cells.Item(132, 10).Value = 290925.66  # J132
$ws.Cells.Item(132, 11).Value = 8536.066500000001  # K132
$ws.Cells.Item(132, 12).Value = 872776.98  # L132
$ws.Cells.Item(132, 13).Value = -6006.066500000001  # M132
$ws.Cells.Item(132, 14).Value = -877836.98  # N132
$ws.Cells.Item(135, 8).Value = 13514489  # H135
$ws.Cells.Item(135, 9).Value = 735.125  # I135
$ws.Cells.Item(135, 10).Value = 100002510  # J135
$ws.Cells.Item(135, 11).Value = 6616.125  # K135
$ws.Cells.Item(135, 12).Value = 900022590  # L135
$ws.Cells.Item(135, 13).Value = -4081.125  # M135
$ws.Cells.Item(135, 14).Value = -900027660  # N135
$ws.Cells.Item(138, 8).Value = 1194.33  # H138
$ws.Cells.Item(138, 9).Value = 573.16364  # I138
$ws.Cells.Item(138, 10).Value = 1953.5333  # J138
$ws.Cells.Item(138, 11).Value = 1719.49092  # K138
$ws.Cells.Item(138, 12).Value = 5860.5999  # L138
$ws.Cells.Item(138, 13).Value = 3420.50908  # M138
$ws.Cells.Item(138, 14).Value = -16140.5999  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1367.72  # H74
$ws.Cells.Item(74, 9).Value = 1233.7561  # I74
$ws.Cells.Item(74, 10).Value = 1978  # J74
$ws.Cells.Item(74, 11).Value = 1233.7561  # K74
$ws.Cells.Item(74, 12).Value = 1978  # L74
$ws.Cells.Item(74, 13).Value = -359.7561000000001  # M74
$ws.Cells.Item(74, 14).Value = -3726  # N74
$ws.Cells.Item(77, 8).Value = 1367.72  # H77
$ws.Cells.Item(77, 9).Value = 1233.7561  # I77
$ws.Cells.Item(77, 10).Value = 1978  # J77
$ws.Cells.Item(77, 11).Value = 6168.780500000001  # K77
$ws.Cells.Item(77, 12).Value = 9890  # L77
$ws.Cells.Item(77, 13).Value = -1800.780500000001  # M77
$ws.Cells.Item(77, 14).Value = -18626  # N77
$ws.Cells.Item(109, 8).Value = 45377  # H109
$ws.Cells.Item(109, 10).Value = 45377  # J109
$ws.Cells.Item(109, 12).Value = 45377  # L109
$ws.Cells.Item(109, 14).Value = -48151  # N109
$ws.Cells.Item(114, 8).Value = 45945  # H114
$ws.Cells.Item(114, 10).Value = 45945  # J114
$ws.Cells.Item(114, 12).Value = 45945  # L114
$ws.Cells.Item(114, 14).Value = -54623  # N114

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1771.9  # H20
$ws.Cells.Item(20, 9).Value = 1262.5  # I20
$ws.Cells.Item(20, 11).Value = 1262.5  # K20
$ws.Cells.Item(20, 13).Value = -1015.5  # M20
$ws.Cells.Item(105, 8).Value = 2784.2173  # H105
$ws.Cells.Item(105, 9).Value = 2021.3636  # I105
$ws.Cells.Item(105, 10).Value = 3483.5  # J105
$ws.Cells.Item(105, 11).Value = 2021.3636  # K105
$ws.Cells.Item(105, 12).Value = 3483.5  # L105
$ws.Cells.Item(105, 13).Value = -274.3635999999999  # M105
$ws.Cells.Item(105, 14).Value = -6977.5  # N105
$ws.Cells.Item(107, 8).Value = 1762.2106  # H107
$ws.Cells.Item(107, 9).Value = 1599.6207  # I107
$ws.Cells.Item(107, 10).Value = 2286.111  # J107
$ws.Cells.Item(107, 11).Value = 1599.6207  # K107
$ws.Cells.Item(107, 12).Value = 2286.111  # L107
$ws.Cells.Item(107, 13).Value = 320.3793000000001  # M107
$ws.Cells.Item(107, 14).Value = -6126.111  # N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3066.27  # H31
$ws.Cells.Item(31, 9).Value = 854.55  # I31
$ws.Cells.Item(31, 10).Value = 3619.2  # J31
$ws.Cells.Item(31, 11).Value = 854.55  # K31
$ws.Cells.Item(31, 12).Value = 3619.2  # L31
$ws.Cells.Item(31, 13).Value = -559.55  # M31
$ws.Cells.Item(31, 14).Value = -4209.2  # N31
$ws.Cells.Item(34, 8).Value = 3066.27  # H34
$ws.Cells.Item(34, 9).Value = 854.55  # I34
$ws.Cells.Item(34, 10).Value = 3619.2  # J34
$ws.Cells.Item(34, 11).Value = 854.55  # K34
$ws.Cells.Item(34, 12).Value = 3619.2  # L34
$ws.Cells.Item(34, 13).Value = -652.55  # M34
$ws.Cells.Item(34, 14).Value = -4023.2  # N34
$ws.Cells.Item(58, 8).Value = 1490.375  # H58
$ws.Cells.Item(58, 9).Value = 1147.6552  # I58
$ws.Cells.Item(58, 10).Value = 2393.9092  # J58
$ws.Cells.Item(58, 11).Value = 1147.6552  # K58
$ws.Cells.Item(58, 12).Value = 2393.9092  # L58
$ws.Cells.Item(58, 13).Value = -944.6551999999999  # M58
$ws.Cells.Item(58, 14).Value = -2799.9092  # N58
$ws.Cells.Item(99, 8).Value = 3271.111  # H99
$ws.Cells.Item(99, 9).Value = 2364.8  # I99
$ws.Cells.Item(99, 10).Value = 4404  # J99
$ws.Cells.Item(99, 11).Value = 2364.8  # K99
$ws.Cells.Item(99, 12).Value = 4404  # L99
$ws.Cells.Item(99, 13).Value = -866.8000000000002  # M99
$ws.Cells.Item(99, 14).Value = -7400  # N99
$ws.Cells.Item(126, 8).Value = 3271.111  # H126
$ws.Cells.Item(126, 9).Value = 2364.8  # I126
$ws.Cells.Item(126, 10).Value = 4404  # J126
$ws.Cells.Item(126, 11).Value = 7094.400000000001  # K126
$ws.Cells.Item(126, 12).Value = 13212  # L126
$ws.Cells.Item(126, 13).Value = -4624.400000000001  # M126
$ws.Cells.Item(126, 14).Value = -18152  # N126
$ws.Cells.Item(132, 8).Value = 41905.37  # H132
$ws.Cells.Item(132, 9).Value = 1569.375  # I132
$ws.Cells.Item(132, 10).Value = 129911.18  # J132
$ws.Cells.Item(132, 11).Value = 4708.125  # K132
$ws.Cells.Item(132, 12).Value = 389733.54  # L132
$ws.Cells.Item(132, 13).Value = -2178.125  # M132
$ws.Cells.Item(132, 14).Value = -394793.54  # N132
$ws.Cells.Item(134, 8).Value = 342682.78  # H134
$ws.Cells.Item(134, 9).Value = 1077.7188  # I134
$ws.Cells.Item(134, 10).Value = 1557278.5  # J134
$ws.Cells.Item(134, 11).Value = 3233.1564  # K134
$ws.Cells.Item(134, 12).Value = 4671835.5  # L134
$ws.Cells.Item(134, 13).Value = -698.1564000000003  # M134
$ws.Cells.Item(134, 14).Value = -4676905.5  # N134
$ws.Cells.Item(136, 8).Value = 1490.375  # H136
$ws.Cells.Item(136, 9).Value = 1147.6552  # I136
$ws.Cells.Item(136, 10).Value = 2393.9092  # J136
$ws.Cells.Item(136, 11).Value = 3442.9656  # K136
$ws.Cells.Item(136, 12).Value = 7181.7276  # L136
$ws.Cells.Item(136, 13).Value = -892.9655999999995  # M136
$ws.Cells.Item(136, 14).Value = -12281.7276  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 5694.048  # H5
$ws.Cells.Item(5, 9).Value = 8824.166999999999  # I5
$ws.Cells.Item(5, 11).Value = 26472.501  # K5
$ws.Cells.Item(5, 13).Value = -26360.501  # M5
$ws.Cells.Item(113, 8).Value = 4632.846  # H113
$ws.Cells.Item(113, 9).Value = 9633.637000000001  # I113
$ws.Cells.Item(113, 10).Value = 965.6  # J113
$ws.Cells.Item(113, 11).Value = 28900.911  # K113
$ws.Cells.Item(113, 12).Value = 2896.8  # L113
$ws.Cells.Item(113, 13).Value = -26730.911  # M113
$ws.Cells.Item(113, 14).Value = -7236.8  # N113
$ws.Cells.Item(135, 8).Value = 5694.048  # H135
$ws.Cells.Item(135, 9).Value = 8824.166999999999  # I135
$ws.Cells.Item(135, 11).Value = 79417.503  # K135
$ws.Cells.Item(135, 13).Value = -76882.503  # M135
$ws.Cells.Item(141, 8).Value = 100002744  # H141
$ws.Cells.Item(141, 9).Value = 111113700  # I141
$ws.Cells.Item(141, 10).Value = 4200  # J141
$ws.Cells.Item(141, 11).Value = 333341100  # K141
$ws.Cells.Item(141, 12).Value = 12600  # L141
$ws.Cells.Item(141, 13).Value = -333335920  # M141
$ws.Cells.Item(141, 14).Value = -22960  # N141

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3878.6  # H40
$ws.Cells.Item(40, 9).Value = 2699.9092  # I40
$ws.Cells.Item(40, 10).Value = 7120  # J40
$ws.Cells.Item(40, 11).Value = 2699.9092  # K40
$ws.Cells.Item(40, 12).Value = 7120  # L40
$ws.Cells.Item(40, 13).Value = -2563.9092  # M40
$ws.Cells.Item(40, 14).Value = -7392  # N40
$ws.Cells.Item(81, 8).Value = 31090.5  # H81
$ws.Cells.Item(81, 10).Value = 31090.5  # J81
$ws.Cells.Item(81, 12).Value = 31090.5  # L81
$ws.Cells.Item(81, 14).Value = -33086.5  # N81
$ws.Cells.Item(84, 8).Value = 31090.5  # H84
$ws.Cells.Item(84, 10).Value = 31090.5  # J84
$ws.Cells.Item(84, 12).Value = 93271.5  # L84
$ws.Cells.Item(84, 14).Value = -103255.5  # N84
$ws.Cells.Item(93, 8).Value = 1513.5  # H93
$ws.Cells.Item(93, 9).Value = 1046.75  # I93
$ws.Cells.Item(93, 10).Value = 1913.5714  # J93
$ws.Cells.Item(93, 11).Value = 1046.75  # K93
$ws.Cells.Item(93, 12).Value = 1913.5714  # L93
$ws.Cells.Item(93, 13).Value = 201.25  # M93
$ws.Cells.Item(93, 14).Value = -4409.5714  # N93
$ws.Cells.Item(101, 8).Value = 25332.334  # H101
$ws.Cells.Item(101, 10).Value = 25332.334  # J101
$ws.Cells.Item(101, 12).Value = 25332.334  # L101
$ws.Cells.Item(101, 14).Value = -31822.334  # N101

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 14663.333  # H70
$ws.Cells.Item(70, 9).Value = 5000  # I70
$ws.Cells.Item(70, 10).Value = 19495  # J70
$ws.Cells.Item(70, 11).Value = 5000  # K70
$ws.Cells.Item(70, 12).Value = 19495  # L70
$ws.Cells.Item(70, 13).Value = -4685  # M70
$ws.Cells.Item(70, 14).Value = -20125  # N70
$ws.Cells.Item(73, 8).Value = 14663.333  # H73
$ws.Cells.Item(73, 9).Value = 5000  # I73
$ws.Cells.Item(73, 10).Value = 19495  # J73
$ws.Cells.Item(73, 11).Value = 5000  # K73
$ws.Cells.Item(73, 12).Value = 19495  # L73
$ws.Cells.Item(73, 13).Value = -3908  # M73
$ws.Cells.Item(73, 14).Value = -21679  # N73
$ws.Cells.Item(132, 8).Value = 1369.7755  # H132
$ws.Cells.Item(132, 9).Value = 1004.09753  # I132
$ws.Cells.Item(132, 10).Value = 3243.875  # J132
$ws.Cells.Item(132, 11).Value = 3012.29259  # K132
$ws.Cells.Item(132, 12).Value = 9731.625  # L132
$ws.Cells.Item(132, 13).Value = -482.29259  # M132
$ws.Cells.Item(132, 14).Value = -14791.625  # N132
$ws.Cells.Item(136, 8).Value = 278521.62  # H136
$ws.Cells.Item(136, 9).Value = 303626.72  # I136
$ws.Cells.Item(136, 10).Value = 2365.6667  # J136
$ws.Cells.Item(136, 11).Value = 910880.1599999999  # K136
$ws.Cells.Item(136, 12).Value = 7097.000100000001  # L136
$ws.Cells.Item(136, 13).Value = -908330.1599999999  # M136
$ws.Cells.Item(136, 14).Value = -12197.0001  # N136
